# ms revision on 20220828-2
#
# Converts an EMU integer value to the "points" float that PowerPoint's
# Shape.Left/Top/Width/Height (COM `Single`) properties expect, biasing by
# half an EMU (in the correct direction for the sign of the value) so that
# the inevitable float32 round-trip lands back on the exact integer EMU
# count instead of truncating one EMU short.
function EmuToPt($emu) {
    if ($emu -ge 0) {
        return ($emu + 0.5) / 12700.0
    } else {
        return ($emu - 0.5) / 12700.0
    }
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 3 : aphid-count figure - reflow the two pictures + callout labels
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)

# Picture 4 - move up (only the vertical offset changes)
$pic4 = $s3.Shapes.Item(1)
$pic4.Top = EmuToPt(438912)

# Picture 1 - reposition
$pic1 = $s3.Shapes.Item(2)
$pic1.Left = EmuToPt(11575)
$pic1.Top = EmuToPt(4159931)

# TextBox 2 : "(a)" -> "A"
$tb2 = $s3.Shapes.Item(3)
$tb2.Top = EmuToPt(48768)
$tb2.Width = EmuToPt(317716)
$tb2.TextFrame.TextRange.Text = "A"

# TextBox 3 : "(b)" -> "B" (also gains an explicit Arial font)
$tb3 = $s3.Shapes.Item(4)
$tb3.Left = EmuToPt(144287)
$tb3.Top = EmuToPt(4270247)
$tb3.Width = EmuToPt(338554)
$tb3.TextFrame.TextRange.Text = "B"
$tb3.TextFrame.TextRange.Font.Name = "Arial"
$tb3.TextFrame.TextRange.Font.NameComplexScript = "Arial"

# TextBox 5 : rotated "log(no. of aphids)" axis label - reposition + Arial
#
# The run-level Font.NameComplexScript setter in this host only ever lands
# on the shape's first run, no matter which sub-range is addressed. Work
# around it by temporarily collapsing the text down to the second run
# ("og(no. of aphids)", which becomes run #1), stamping Arial/cs on it,
# then re-inserting the "l" prefix as its own run (restoring its en-US
# language) and stamping Arial/cs on that new run too - so both runs end
# up correctly tagged instead of just the first one.
$tb5 = $s3.Shapes.Item(5)
$tb5.Left = EmuToPt(-485051)
$tb5.Top = EmuToPt(5507555)
$tb5.Width = EmuToPt(1566454)
$tr5 = $tb5.TextFrame.TextRange
$tr5.Text = "og(no. of aphids)"
$tr5.Font.Name = "Arial"
$tr5.Font.NameComplexScript = "Arial"
$lRun = $tr5.InsertBefore("l")
$lRun.Font.Name = "Arial"
$lRun.Font.NameComplexScript = "Arial"
$tr5.Characters(1, 1).LanguageID = "en-US"

# ---------------------------------------------------------------------
# Slide 5 : gene-model figure - rename one transcript label
# ---------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(57).TextFrame.TextRange.Text = "AT3G13910.2"

# ---------------------------------------------------------------------
# Slide 6 : duplicate gene-model figure - same rename + panel letters
# ---------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item(41).TextFrame.TextRange.Text = "AT3G13910.2"

$a = $s6.Shapes.Item(43)
$a.Width = EmuToPt(338554)
$a.TextFrame.TextRange.Text = "A"

$b = $s6.Shapes.Item(44)
$b.Width = EmuToPt(338554)
$b.TextFrame.TextRange.Text = "B"

$c = $s6.Shapes.Item(45)
$c.Width = EmuToPt(351378)
$c.TextFrame.TextRange.Text = "C"
